$wb = $excel.ActiveWorkbook

# --- 1. Rename sheet AddProductCategory1 -> AddProductCategory ---
$wsP = $wb.Worksheets.Item("AddProductCategory1")
$wsP.Name = "AddProductCategory"

# Fix the stale RefersTo on the pre-existing Print_Area defined name left over
# from before the rename (renaming a sheet does not rewrite old formula text).
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area") {
        $n.RefersTo = "=AddProductCategory!`$13:`$13"
    }
}

# Re-apply the print area a few more times (mirrors the extra duplicate
# _xlnm.Print_Area / _xlnm.Print_Area_0 / _xlnm.Print_Area_0_0 entries seen
# in the saved workbook).
$wsP.Names.Add("_xlnm.Print_Area", "=AddProductCategory!`$13:`$13")
$wsP.Names.Add("_xlnm.Print_Area_0", "=AddProductCategory!`$13:`$13")
$wsP.Names.Add("_xlnm.Print_Area_0_0", "=AddProductCategory!`$13:`$13")

# --- 2. "Web Data 1" -> "Web Data 11" everywhere it appears ---
foreach ($ws in $wb.Worksheets) {
    $c = $ws.Range("C1")
    if ($c.Text -eq "Web Data 1") {
        $c.Value = "Web Data 11"
    }
}

# --- 3. AddCustomer: rename dependent-holder names, drop the Jayden row ---
$wsC = $wb.Worksheets.Item("AddCustomer")
$wsC.Range("D1").Value = "LEO DEPENDENT"
$wsC.Range("E1").Value = "Leo Dependent"
$wsC.Range("D2").Value = "WILLIAM DEPENDENT"
$wsC.Range("E2").Value = "William Dependent"
$wsC.Range("D3").Value = "OLIVIA DEPENDENT"
$wsC.Range("E3").Value = "Olivia Dependent"
$wsC.Range("D4").Value = "ISABELLA DEPENDENT"
$wsC.Range("E4").Value = "Isabella Dependent"
$wsC.Range("D5").Value = "SOPHIA DEPENDENT"
$wsC.Range("E5").Value = "Sophia Dependent"
$wsC.Range("D6").Value = "JACOB DEPENDENT"
$wsC.Range("E6").Value = "Jacob Dependent"

# Row 4 grew slightly taller
$wsC.Rows.Item(4).RowHeight = 16.25

# The old row 7 (Jayden) is removed entirely
$wsC.Rows.Item(7).Delete()

# Selection moves from E18 to E13
$wsC.Range("E13").Select()

# --- 4. CreateOrder: same dependent-holder renames (normal-case column) ---
$wsO = $wb.Worksheets.Item("CreateOrder")
$wsO.Range("D1").Value = "Leo Dependent"
$wsO.Range("D2").Value = "William Dependent"
$wsO.Range("D3").Value = "Olivia Dependent"
$wsO.Range("D4").Value = "Isabella Dependent"
$wsO.Range("D5").Value = "Sophia Dependent"
$wsO.Range("D6").Value = "Jacob Dependent"

# Row 4 grew slightly taller (same as AddCustomer)
$wsO.Rows.Item(4).RowHeight = 16.25

# Column D widened
$wsO.Columns.Item(4).ColumnWidth = 20.6814814814815

# Selection / top-left cell moves from C1 to D1
$wsO.Range("D1").Select()
$excel.ActiveWindow.ScrollColumn = 4

# --- 5. Re-select AddProductCategory last so it stays the active tab ---
# (selection there moves from A8 to C1)
$wsP.Activate()
$wsP.Range("C1").Select()
